$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts batsman..sr to F..K),
# making room for the new "ownTeam" and "oppTeam" columns.
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(4).Insert()

# New headers for the inserted columns.
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# Keep the numeric-looking columns (totalRuns, totalBalls, total4s, total6s, sr)
# stored as text, matching the original data which used text cells for every
# value (t="str"), so that values like "90.00" or "0" keep their exact text.
$ws.Range("G2:K15").NumberFormat = "@"

$data = New-Object 'object[,]' 14,11
$data[0,0] = ' Abu Dhabi'
$data[0,1] = ' October 07 2020'
$data[0,2] = 'KKR won by 10 runs'
$data[0,3] = 'Kolkata Knight Riders'
$data[0,4] = 'Chennai Super Kings'
$data[0,5] = 'Nitish Rana '
$data[0,6] = '9'
$data[0,7] = '10'
$data[0,8] = '1'
$data[0,9] = '0'
$data[0,10] = '90.00'
$data[1,0] = ' Dubai (DSC)'
$data[1,1] = ' September 30 2020'
$data[1,2] = 'KKR won by 37 runs'
$data[1,3] = 'Kolkata Knight Riders'
$data[1,4] = 'Rajasthan Royals'
$data[1,5] = 'Nitish Rana '
$data[1,6] = '22'
$data[1,7] = '17'
$data[1,8] = '2'
$data[1,9] = '1'
$data[1,10] = '129.41'
$data[2,0] = ' Abu Dhabi'
$data[2,1] = ' October 16 2020'
$data[2,2] = 'Mumbai won by 8 wickets (with 19 balls remaining)'
$data[2,3] = 'Kolkata Knight Riders'
$data[2,4] = 'Mumbai Indians'
$data[2,5] = 'Nitish Rana '
$data[2,6] = '5'
$data[2,7] = '6'
$data[2,8] = '1'
$data[2,9] = '0'
$data[2,10] = '83.33'
$data[3,0] = ' Abu Dhabi'
$data[3,1] = ' October 18 2020'
$data[3,2] = 'Match tied (KKR won the one-over eliminator)'
$data[3,3] = 'Kolkata Knight Riders'
$data[3,4] = 'Sunrisers Hyderabad'
$data[3,5] = 'Nitish Rana '
$data[3,6] = '29'
$data[3,7] = '20'
$data[3,8] = '3'
$data[3,9] = '1'
$data[3,10] = '145.00'
$data[4,0] = ' Abu Dhabi'
$data[4,1] = ' October 10 2020'
$data[4,2] = 'KKR won by 2 runs'
$data[4,3] = 'Kolkata Knight Riders'
$data[4,4] = 'Kings XI Punjab'
$data[4,5] = 'Nitish Rana '
$data[4,6] = '2'
$data[4,7] = '4'
$data[4,8] = '0'
$data[4,9] = '0'
$data[4,10] = '50.00'
$data[5,0] = ' Dubai (DSC)'
$data[5,1] = ' November 01 2020'
$data[5,2] = 'KKR won by 60 runs'
$data[5,3] = 'Kolkata Knight Riders'
$data[5,4] = 'Rajasthan Royals'
$data[5,5] = 'Nitish Rana '
$data[5,6] = '0'
$data[5,7] = '1'
$data[5,8] = '0'
$data[5,9] = '0'
$data[5,10] = '0.00'
$data[6,0] = ' Abu Dhabi'
$data[6,1] = ' September 26 2020'
$data[6,2] = 'KKR won by 7 wickets (with 12 balls remaining)'
$data[6,3] = 'Kolkata Knight Riders'
$data[6,4] = 'Sunrisers Hyderabad'
$data[6,5] = 'Nitish Rana '
$data[6,6] = '26'
$data[6,7] = '13'
$data[6,8] = '6'
$data[6,9] = '0'
$data[6,10] = '200.00'
$data[7,0] = ' Sharjah'
$data[7,1] = ' October 03 2020'
$data[7,2] = 'Capitals won by 18 runs'
$data[7,3] = 'Kolkata Knight Riders'
$data[7,4] = 'Delhi Capitals'
$data[7,5] = 'Nitish Rana '
$data[7,6] = '58'
$data[7,7] = '35'
$data[7,8] = '4'
$data[7,9] = '4'
$data[7,10] = '165.71'
$data[8,0] = ' Dubai (DSC)'
$data[8,1] = ' October 29 2020'
$data[8,2] = 'Super Kings won by 6 wickets'
$data[8,3] = 'Kolkata Knight Riders'
$data[8,4] = 'Chennai Super Kings'
$data[8,5] = 'Nitish Rana '
$data[8,6] = '87'
$data[8,7] = '61'
$data[8,8] = '10'
$data[8,9] = '4'
$data[8,10] = '142.62'
$data[9,0] = ' Sharjah'
$data[9,1] = ' October 26 2020'
$data[9,2] = 'Kings XI won by 8 wickets (with 7 balls remaining)'
$data[9,3] = 'Kolkata Knight Riders'
$data[9,4] = 'Kings XI Punjab'
$data[9,5] = 'Nitish Rana '
$data[9,6] = '0'
$data[9,7] = '1'
$data[9,8] = '0'
$data[9,9] = '0'
$data[9,10] = '0.00'
$data[10,0] = ' Abu Dhabi'
$data[10,1] = ' October 21 2020'
$data[10,2] = 'RCB won by 8 wickets (with 39 balls remaining)'
$data[10,3] = 'Kolkata Knight Riders'
$data[10,4] = 'Royal Challengers Bangalore'
$data[10,5] = 'Nitish Rana '
$data[10,6] = '0'
$data[10,7] = '1'
$data[10,8] = '0'
$data[10,9] = '0'
$data[10,10] = '0.00'
$data[11,0] = ' Abu Dhabi'
$data[11,1] = ' October 24 2020'
$data[11,2] = 'KKR won by 59 runs'
$data[11,3] = 'Kolkata Knight Riders'
$data[11,4] = 'Delhi Capitals'
$data[11,5] = 'Nitish Rana '
$data[11,6] = '81'
$data[11,7] = '53'
$data[11,8] = '13'
$data[11,9] = '1'
$data[11,10] = '152.83'
$data[12,0] = ' Abu Dhabi'
$data[12,1] = ' September 23 2020'
$data[12,2] = 'Mumbai won by 49 runs'
$data[12,3] = 'Kolkata Knight Riders'
$data[12,4] = 'Mumbai Indians'
$data[12,5] = 'Nitish Rana '
$data[12,6] = '24'
$data[12,7] = '18'
$data[12,8] = '2'
$data[12,9] = '1'
$data[12,10] = '133.33'
$data[13,0] = ' Sharjah'
$data[13,1] = ' October 12 2020'
$data[13,2] = 'RCB won by 82 runs'
$data[13,3] = 'Kolkata Knight Riders'
$data[13,4] = 'Royal Challengers Bangalore'
$data[13,5] = 'Nitish Rana '
$data[13,6] = '9'
$data[13,7] = '14'
$data[13,8] = '1'
$data[13,9] = '0'
$data[13,10] = '64.28'

$ws.Range("A2:K15").Value = $data
